$wb = $excel.ActiveWorkbook

$login = $wb.Worksheets.Item("Login")
$recordings = $wb.Worksheets.Item("Recordings")

# --- Create the new "Logout" sheet as a copy of "Recordings", placed last ---
$recordings.Copy($null, $recordings)
$logout = $wb.Worksheets.Item($wb.Worksheets.Count)
$logout.Name = "Logout"

# Drop the extra columns G:I that "Recordings" has but "Logout" should not.
$logout.Range("G1:I2").EntireColumn.Delete()

# Overwrite row 2 with the same scenario data used on the "Login" sheet.
$logout.Range("A2").Value = $login.Range("A2").Value2
$logout.Range("B2").Value = $login.Range("B2").Value2
$logout.Range("C2").Value = $login.Range("C2").Value2
$logout.Range("D2").Value = $login.Range("D2").Value2
$logout.Range("E2").Value = $login.Range("E2").Value2
$logout.Range("F2").Value = $login.Range("F2").Value2

# --- Update selections on the other two sheets ---
$login.Range("A1").Select()
$recordings.Range("G1").Select()

# --- Leave "Logout" as the active sheet/selection ---
$logout.Range("F2").Select()
